# Update countries & provincias Spain
#
# 1) Reorder "Nueva Caledonia" / "Santa Lucia" in the country list
#    (Nueva Caledonia now appears before Santa Lucia).
# 2) Refresh the "Datos actualizados..." timestamp (04:40 -> 05:57).
# 3) Update the day's case numbers for India, Belgica, Kazajistan,
#    Mongolia and Butan.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap Santa Lucia (row 207) and Nueva Caledonia (row 208) ---
$ws.Range("A207").Value = "Nueva Caledonia"
$ws.Range("A208").Value = "Santa Lucia"

# --- 2) Update "last refreshed" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Octubre de 2020 a las 05:57"

# --- 3) Update case counters ---

# Row 5: India
$ws.Range("B5").Value = 6757131
$ws.Range("C5").Value = 2952
$ws.Range("D5").Value = 5744693
$ws.Range("E5").Value = 907847

# Row 35: Belgica
$ws.Range("B35").Value = 134291
$ws.Range("C35").Value = 2088
$ws.Range("D35").Value = 19812
$ws.Range("E35").Value = 104387
$ws.Range("G35").Value = 14
$ws.Range("H35").Value = 10092

# Row 39: Kazajistan
$ws.Range("D39").Value = 103465
$ws.Range("E39").Value = 3151

# Row 186: Mongolia
$ws.Range("D186").Value = 308
$ws.Range("E186").Value = 7

# Row 187: Butan
$ws.Range("B187").Value = 300
$ws.Range("C187").Value = 1
$ws.Range("D187").Value = 250
$ws.Range("E187").Value = 50
